# Ajout slide sur la vérification de code
#
# The "F. Résultats et analyse - Avant" slide (slide 8) is duplicated.
# The duplicate takes slide 8's old spot (position 8) and is rewritten
# into the new "G. Vérification de code" slide (title + a single text
# box, no pictures). The original slide is moved to the end of the
# deck (now position 10).

$p = $ppt.ActivePresentation

$original = $p.Slides.Item(8)
$newSlide = $original.Duplicate().Item(1)

# Move the original slide (still holding the old "F." content) to the
# very end of the deck; this shifts $newSlide up into position 8.
$original.MoveTo($p.Slides.Count)

# --- Rewrite the duplicate into the new "G. Vérification de code" slide ---

# Title
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "G. Vérification de code"

# Remove the two picture shapes copied from the original slide (by
# name, so shape order doesn't matter).
for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $newSlide.Shapes.Item($i)
    if ($shp.Name -eq "Picture 14" -or $shp.Name -eq "Image 5") {
        $shp.Delete()
    }
}

# Body text box: replace its text/formatting, then reposition/resize
# it last (the box has spAutoFit, which recomputes Height as the text
# changes, so Height must be forced after the text is final).
$body = $newSlide.Shapes.Item("ZoneTexte 8")

$bodyText = "Un module de test a été ajouté au projet, à l’aide du framework python unittest, afin de pouvoir ajouter un ensemble de test facilement pour pouvoir tester le code. Pour ce qui est de la vérification, les ordres de précision asymptotiques sont un bon indicateur pour vérifier le code. Additionnellement, un test d’invariance galiléenne a été ajouté pour s’assurer que le code fonctionne correctement. D’autres tests de profilage ont été réalisés afin de s’assurer que le code produit des résultats valide (non nuls ou None) en sortie."

$tr = $body.TextFrame.TextRange
$tr.Text = $bodyText
$tr.Font.Size = 20

# "framework" is its own run in the source deck (flagged by the spell
# checker) even though its visible formatting matches its neighbours;
# force a run break there by re-asserting the (unchanged) font size.
$beforeFramework = "Un module de test a été ajouté au projet, à l’aide du "
$frameworkStart = $beforeFramework.Length + 1
$body.TextFrame.TextRange.Characters($frameworkStart, ("framework").Length).Font.Size = 20

# "unittest" -> italic (and its own run, flagged by the spell checker)
$beforeUnittest = "Un module de test a été ajouté au projet, à l’aide du framework python "
$unittestStart = $beforeUnittest.Length + 1
$body.TextFrame.TextRange.Characters($unittestStart, ("unittest").Length).Font.Italic = $true

# The closing sentence ("ce qui est de la vérification...") is its own
# run too (it keeps the inherited noProof flag from the original
# slide); force the run break the same way.
$tailText = "ce qui est de la vérification, les ordres de précision asymptotiques sont un bon indicateur pour vérifier le code. Additionnellement, un test d’invariance galiléenne a été ajouté pour s’assurer que le code fonctionne correctement. D’autres tests de profilage ont été réalisés afin de s’assurer que le code produit des résultats valide (non nuls ou None) en sortie."
$tailStart = $bodyText.Length - $tailText.Length + 1
$body.TextFrame.TextRange.Characters($tailStart, $tailText.Length).Font.Size = 20

# Final position/size (EMU 833674, 2175263, 10520126, 1938992 -> pt).
$body.Left = 65.64366929133858
$body.Top = 171.28059842519684
$body.Width = 828.3564251968504
$body.Height = 152.67658267716533
